$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add I7 and I9 values (new data points)
$ws.Range("I7").Value = 0.5
$ws.Range("I9").Value = 0.5

# Update H10 value
$ws.Range("H10").Value = 4.5

# Add I10 and J10 values
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 3

# Add per-row total formulas in column B for rows 3..10
$ws.Range("B3").Formula = "=SUM(C3:S3)"
$ws.Range("B4:B10").Formula = "=SUM(C4:S4)"

# Add grand total formula in B11
$ws.Range("B11").Formula = "=SUM(B3:B10)"

# Update selection to J10
$ws.Range("J10").Select()
